$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove one blank template row so the 6 leftover blank rows (11-16)
#    become 5, matching the 5 new diary entries added beneath the existing
#    one (dimension goes from A1:G125 -> A1:G124).
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# 2. Propagate the filled "diary entry" row formatting (row 10) down onto
#    the five rows that will hold the newly written diary entries.
# ---------------------------------------------------------------------------
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Update the Time for the first (existing) entry.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "05:00 pm - 07:50 pm"
$ws.Range("C10").Value = "N/A"

# ---------------------------------------------------------------------------
# 4. Fill in the five new diary entries.
# ---------------------------------------------------------------------------

# Row 11 - 2020-01-13
$ws.Range("A11").Value = 43843
$ws.Range("B11").Value = "10:00 pm - 10:30 pm"
$ws.Range("C11").Value = "Anjana, Aman"
$ws.Range("D11").Value = "Team formation"
$ws.Range("E11").Value = "Formed a group of three members for the project and added the team information on Excel sheet"
$ws.Range("F11").Value = "Need to learn git properly. "
$ws.Range("G11").Value = "Excited to be working with my team for the rest of the quarter"

# Row 12 - 2020-01-14
$ws.Range("A12").Value = 43844
$ws.Range("B12").Value = "11:00 am - 02:00 pm"
$ws.Range("C12").Value = "N/A"
$ws.Range("D12").Value = "To successfully build and run 3 projects."
$ws.Range("E12").Value = "Was able to run 2 projects. "
$ws.Range("F12").Value = "Had some trouble building the projects. Maybe I should've read the readme thoroughly first."
$ws.Range("G12").Value = "Little frustated about fixing the errors. Errors were a little random."

# Row 13 - 2020-01-16
$ws.Range("A13").Value = 43846
$ws.Range("B13").Value = "05:00 pm - 07:50 pm"
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "Didn't know what to expect. Probably learn more about code reading and understanding."
$ws.Range("E13").Value = "Was able to check the code to understand its functionality briefly by successfully navigating through the project. Felt motivated after Mr.Ping interaction with the class."
$ws.Range("F13").Value = "Search is a good tool to navigate the  files efficiently. Also, find usages was effective."
$ws.Range("G13").Value = "Motivated! "

# Row 14 - 2020-01-18
$ws.Range("A14").Value = 43848
$ws.Range("B14").Value = "06:00 pm - 10:00 pm"
$ws.Range("C14").Value = "Anjana, Aman"
$ws.Range("D14").Value = "Select an open source project on GitHub"
$ws.Range("E14").Value = "Explored various open source projects available on GitHub, selected OpenRefine and submitted the pull request. "
$ws.Range("F14").Value = "Narrowing down one project from a list of possible options was more difficult that expected. "
$ws.Range("G14").Value = "A little unsure of the selected project."

# Row 15 - 2020-01-19
$ws.Range("A15").Value = 43849
$ws.Range("B15").Value = "11:00 am - 12:30 pm"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "To complete the individual homework assigned. "
$ws.Range("E15").Value = "Completed the given homework."
$ws.Range("F15").Value = "Figured out some details for adding the fruits to the game, though unsure of some details."
$ws.Range("G15").Value = "Hopefully what I submitted was correct."

# ---------------------------------------------------------------------------
# 5. Row heights - wrap-text autofit heights as they appear in the saved
#    workbook.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 93.6
$ws.Rows.Item(11).RowHeight = 46.8
$ws.Rows.Item(12).RowHeight = 46.8
$ws.Rows.Item(13).RowHeight = 78
$ws.Rows.Item(14).RowHeight = 62.4
$ws.Rows.Item(15).RowHeight = 46.8

# ---------------------------------------------------------------------------
# 6. View state - scrolled/selected cell as last left by the author.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D16").Select()
